# Fruta / hortaliza, semanal
# Insert two new weekly price rows for "Cereza" (Lapins / Rainier, Macroferia
# Regional de Talca) above the existing row 111, pushing the rest of the
# table down by two rows (171 -> 173 total data rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above what is currently row 111. Doing this twice at
# the same index pushes all the old rows (111..171) down to (113..173).
$ws.Rows.Item(111).Insert()
$ws.Rows.Item(111).Insert()

# --- New row 111: Cereza / Lapins / Primera ---
$ws.Cells.Item(111, 1).Value = 5
$ws.Cells.Item(111, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(111, 3).Value = "Maule"
$ws.Cells.Item(111, 4).Value = 44572
$ws.Cells.Item(111, 5).Value = 7
$ws.Cells.Item(111, 6).Value = "Fruta"
$ws.Cells.Item(111, 7).Value = 100103
$ws.Cells.Item(111, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(111, 9).Value = 100103001
$ws.Cells.Item(111, 10).Value = "Cereza"
$ws.Cells.Item(111, 11).Value = "Lapins"
$ws.Cells.Item(111, 12).Value = "Primera"
$ws.Cells.Item(111, 13).Value = 440
$ws.Cells.Item(111, 14).Value = 4000
$ws.Cells.Item(111, 15).Value = 5000
$ws.Cells.Item(111, 16).Value = 4545
$ws.Cells.Item(111, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(111, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(111, 19).Value = 454
$ws.Cells.Item(111, 20).Value = 10

# --- New row 112: Cereza / Rainier / Primera ---
$ws.Cells.Item(112, 1).Value = 5
$ws.Cells.Item(112, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(112, 3).Value = "Maule"
$ws.Cells.Item(112, 4).Value = 44572
$ws.Cells.Item(112, 5).Value = 7
$ws.Cells.Item(112, 6).Value = "Fruta"
$ws.Cells.Item(112, 7).Value = 100103
$ws.Cells.Item(112, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(112, 9).Value = 100103001
$ws.Cells.Item(112, 10).Value = "Cereza"
$ws.Cells.Item(112, 11).Value = "Rainier"
$ws.Cells.Item(112, 12).Value = "Primera"
$ws.Cells.Item(112, 13).Value = 160
$ws.Cells.Item(112, 14).Value = 5500
$ws.Cells.Item(112, 15).Value = 5500
$ws.Cells.Item(112, 16).Value = 5500
$ws.Cells.Item(112, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(112, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(112, 19).Value = 550
$ws.Cells.Item(112, 20).Value = 10
